$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update confidential notice date from 2021-07-08 to 2021-07-09
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-09 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) columns for holdings rows 2-56
$ws.Range("D2").Value = 0.01874460964046005
$ws.Range("E2").Value = 0.03197320341047494
$ws.Range("D3").Value = 0.0181281608459504
$ws.Range("E3").Value = 0.04617968094038627
$ws.Range("D4").Value = 0.02025719233069205
$ws.Range("E4").Value = 0.0263736263736265
$ws.Range("D5").Value = 0.01990901291897828
$ws.Range("E5").Value = 0.02614678899082556
$ws.Range("D6").Value = 0.01997979037316272
$ws.Range("E6").Value = 0.03599588618443605
$ws.Range("D7").Value = 0.008027152666778415
$ws.Range("E7").Value = 0.008058781701825168
$ws.Range("D8").Value = 0.01960459376119573
$ws.Range("E8").Value = 0.02477678571428577
$ws.Range("D9").Value = 0.02327094199273922
$ws.Range("E9").Value = 0.0264900662251657
$ws.Range("D10").Value = 0.0229775580294263
$ws.Range("E10").Value = 0.02027027027027017
$ws.Range("D11").Value = 0.01921836195475913
$ws.Range("E11").Value = 0.05280665280665287
$ws.Range("D12").Value = 0.014368584247336
$ws.Range("E12").Value = 0.006091101694915446
$ws.Range("D13").Value = 0.01495230798238402
$ws.Range("E13").Value = 0.009619789280806312
$ws.Range("D14").Value = 0.009020700692991188
$ws.Range("E14").Value = 0.003712140386399998
$ws.Range("D15").Value = 0.0143929377799586
$ws.Range("E15").Value = 0.009835025380710682
$ws.Range("D16").Value = 0.02374374299717024
$ws.Range("E16").Value = 0.03129131776112826
$ws.Range("D17").Value = 0.02479341830559917
$ws.Range("E17").Value = 0.0201976794155565
$ws.Range("D18").Value = 0.02298402693652918
$ws.Range("E18").Value = 0.03208556149732633
$ws.Range("D19").Value = 0.01759904229732962
$ws.Range("E19").Value = 0.02497324295397796
$ws.Range("D20").Value = 0.01962209786276823
$ws.Range("E20").Value = 0.0373986735445837
$ws.Range("D21").Value = 0.02820443496855255
$ws.Range("E21").Value = 0.003669724770642091
$ws.Range("D22").Value = 0.01832907749008688
$ws.Range("E22").Value = 0.007006726457399193
$ws.Range("D23").Value = 0.02019592797518832
$ws.Range("E23").Value = 0.04299657082563946
$ws.Range("D24").Value = 0.0185697588867087
$ws.Range("E24").Value = 0.03814510097232615
$ws.Range("D25").Value = 0.01986715528478317
$ws.Range("E25").Value = -0.006607929515418554
$ws.Range("D26").Value = 0.01812530691634619
$ws.Range("E26").Value = 0.01684774051330495
$ws.Range("D27").Value = 0.02247831061066293
$ws.Range("E27").Value = 0.0255281690140845
$ws.Range("D28").Value = 0.02226635877205684
$ws.Range("E28").Value = 0.01999487310945924
$ws.Range("D29").Value = 0.01955797957766028
$ws.Range("E29").Value = 0.02621722846441932
$ws.Range("D30").Value = 0.02164819761978466
$ws.Range("E30").Value = 0.03122665471387998
$ws.Range("D31").Value = 0.01998721059013366
$ws.Range("E31").Value = 0.05462108880448557
$ws.Range("D32").Value = 0.02079677528786162
$ws.Range("E32").Value = 0.007081038552320784
$ws.Range("D33").Value = 0.01758096740983628
$ws.Range("E33").Value = 0.02054023635340485
$ws.Range("D34").Value = 0.018649478653653
$ws.Range("E34").Value = 0.02396449704142012
$ws.Range("D35").Value = 0.02152776179048694
$ws.Range("E35").Value = -0.0004242231413723596
$ws.Range("D36").Value = 0.01776571178621556
$ws.Range("E36").Value = 0.01590361445783128
$ws.Range("D37").Value = 0.0204725688848232
$ws.Range("E37").Value = 0.01523205888366386
$ws.Range("D38").Value = 0.01928971019486441
$ws.Range("E38").Value = 0.01597869507323568
$ws.Range("D39").Value = 0.01878684779860238
$ws.Range("E39").Value = 0.03275202041684389
$ws.Range("D40").Value = 0.0167776813572376
$ws.Range("E40").Value = 0.009797917942437229
$ws.Range("D41").Value = 0.01333508120666427
$ws.Range("E41").Value = 0.02157287980824107
$ws.Range("D42").Value = 0.01454971364621661
$ws.Range("E42").Value = 0.0517836593785963
$ws.Range("D43").Value = 0.01659598117243614
$ws.Range("E43").Value = 0.03683492496589369
$ws.Range("D44").Value = 0.01282955514277163
$ws.Range("E44").Value = 0.02544823597455181
$ws.Range("D45").Value = 0.01510508834719614
$ws.Range("E45").Value = 0.03279968762202268
$ws.Range("D46").Value = 0.02028573162673417
$ws.Range("E46").Value = -0.04351903957981618
$ws.Range("D47").Value = 0.01354608173540229
$ws.Range("E47").Value = 0.04198210598761198
$ws.Range("D48").Value = 0.01985764218610247
$ws.Range("E48").Value = 0.03679218166139675
$ws.Range("D49").Value = 0.01814851887712711
$ws.Range("E49").Value = 0.03368383532347186
$ws.Range("D50").Value = 0.01792305843839441
$ws.Range("E50").Value = 0.002133712660028486
$ws.Range("D51").Value = 0.01916299572043743
$ws.Range("E51").Value = 0.03109641676347064
$ws.Range("D52").Value = 0.006479561773401415
$ws.Range("E52").Value = 0.006342494714587588
$ws.Range("D53").Value = 0.02142939635012846
$ws.Range("E53").Value = 0.01621223286661744
$ws.Range("D54").Value = 0.01761217037350899
$ws.Range("E54").Value = 0.03853383458646609
$ws.Range("D55").Value = 0.02066796793172487
$ws.Range("E55").Value = 0.04107558755028595
$ws.Range("D56").Value = 1
$ws.Range("E56").Value = 0.02338205498533008

$ws.Protect()
